$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 22 for the "diaryEntry" (Tagebucheintrag) sub-entry
# that belongs to the "Art der Einklebung" item in row 21.
$ws.Rows.Item(22).Insert()

# Copy the formatting pattern used by similar "sub-row" entries (row 9 has the
# same A/B empty + C "type" + D description layout) onto the new row.
$ws.Range("A9:D9").Copy()
$ws.Range("A22:D22").PasteSpecial(-4122)

# Row 21 ("Art der Einklebung") no longer carries the old "n" value in column C,
# since that information now lives in the new sub-row below it.
$ws.Range("C21").Value = ""

# Fill in the new sub-row describing the "diaryEntry" markup for diary entries.
$ws.Range("C22").Value = "type"
$ws.Range("D22").Value = "diaryEntry für Tagebucheintrag"

# The data table used by the AutoFilter grew by one row (it now ends at row 67
# instead of row 66), so refresh the filter range accordingly.
$ws.AutoFilterMode = $false
$ws.Range("A4:D67").AutoFilter()

# Keep the workbook-level _FilterDatabase defined name in sync with the new
# AutoFilter range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Tabelle1!_FilterDatabase") {
        $n.RefersTo = "=Tabelle1!`$A`$4:`$D`$67"
    }
}

# Reflect where the editor ended up after making the change.
$ws.Range("D22").Select()
